# Fills in the rest of the weekly workout table (Tuesday..Saturday columns,
# second exercise row, and the trailing "Abs" note) and widens/adds the
# supporting columns so the new text fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

function Set-Cell($addr, $text) {
    $ws.Range($addr).Value = $text
    $ws.Range($addr).HorizontalAlignment = $xlCenter
}

# Row 5 - day headers
Set-Cell "G5" "Tuesday"
Set-Cell "K5" "Wednesday"
Set-Cell "O5" "Thursday"
Set-Cell "S5" "Friday"
Set-Cell "W5" "Saturday"

# Row 7 - workout focus per day
Set-Cell "G7" "Back"
Set-Cell "K7" "Shoulder"
Set-Cell "O7" "Chest & Triceps"
Set-Cell "S7" "Back & Biceps"
Set-Cell "W7" "Legs & Shoulders"

# Row 9 - first exercise, updated sets/reps for Monday, new for Wednesday
Set-Cell "D9" "4 sets"
Set-Cell "E9" "8-12 reps"
Set-Cell "K9" "Box style"
Set-Cell "L9" "4 sets"
Set-Cell "M9" "8-12 reps"

# Row 10 - second exercise, new for Monday and Wednesday
Set-Cell "C10" "Incline Bench Press"
Set-Cell "D10" "3 sets"
Set-Cell "E10" "4-6 reps"
Set-Cell "K10" "Jog"
Set-Cell "L10" "3 sets"
Set-Cell "M10" "4-6 reps"

# Row 12 - trailing note under Saturday
Set-Cell "W12" "Abs"

# Widen / add the columns needed for the new content (values derived from
# Excel's ColumnWidth<->stored-width pixel rounding so the saved widths line
# up with the target worksheet).
$ws.Columns.Item(3).ColumnWidth = 21          # C
$ws.Columns.Item(5).ColumnWidth = 9.5         # E
$ws.Columns.Item(7).ColumnWidth = 7.1666667   # G
$ws.Columns.Item(11).ColumnWidth = 9.5        # K
$ws.Columns.Item(12).ColumnWidth = 6          # L
$ws.Columns.Item(13).ColumnWidth = 9.5        # M
$ws.Columns.Item(15).ColumnWidth = 16.3333333 # O
$ws.Columns.Item(19).ColumnWidth = 14.1666667 # S
$ws.Columns.Item(23).ColumnWidth = 17.5       # W
